# Regenerate the "K" (strikeouts) column (G) of the save-data sheet.
# The previous values were derived from an approximate "Strike#" count;
# this regenerates them to use the actual "K" (strikeouts) stat instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value for column G (rows 2..33)
$kValues = [ordered]@{
    2  = 2
    3  = 4
    4  = 2
    5  = 2
    6  = 2
    7  = 3
    8  = 3
    9  = 4
    10 = 5
    11 = 1
    12 = 2
    13 = 4
    14 = 6
    15 = 2
    16 = 7
    17 = 4
    18 = 6
    19 = 5
    20 = 4
    21 = 4
    22 = 1
    23 = 4
    24 = 5
    25 = 7
    26 = 5
    27 = 2
    28 = 5
    29 = 5
    30 = 7
    31 = 7
    32 = 2
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
